$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.026200652122498
$ws.Range("B1").Value = 3.522655487060547
$ws.Range("C1").Value = 3.856557369232178
$ws.Range("D1").Value = 3.164623975753784
$ws.Range("E1").Value = 1.290530562400818
